# Insert a new data row at row 514 (pushing the existing rows 514-553 down to 515-554),
# mirroring the weekly refresh of the "Ciboulette" price series for
# "Hortaliza, Vega Central Mapocho de Santiago" with one additional observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 514:553 down to 515:554 and leave a blank row 514 behind.
$ws.Rows.Item(514).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A514").Value = 9
$ws.Range("B514").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C514").Value = "Metropolitana"
$ws.Range("D514").Value = 45013
$ws.Range("E514").Value = 13
$ws.Range("F514").Value = 100112039
$ws.Range("G514").Value = "Ciboulette"
$ws.Range("H514").Value = "Sin especificar"
$ws.Range("I514").Value = "Primera"
$ws.Range("J514").Value = 340
$ws.Range("K514").Value = 1000
$ws.Range("L514").Value = 1200
$ws.Range("M514").Value = 1100
$ws.Range("N514").Value = "$/docena de atados"
$ws.Range("O514").Value = "Región Metropolitana"
$ws.Range("P514").Value = 367
$ws.Range("Q514").Value = 3
$ws.Range("R514").Value = "Hortaliza"
